$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge split "UML (Unified Modeling Language)" runs (spell-check
#    artifacts from <w:proofErr> wrapped sub-words) into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "UML (Unified Modeling Language): Linguagem de modelagem que define representações de um sistema de forma padronizada com o objetivo de facilitar a compreensão.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "UML (Unified Modeling Language): Linguagem de modelagem que define representações de um sistema de forma padronizada com o objetivo de facilitar a compreensão.",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) Merge split "GUEDES, Gilleanes T. A." run and the ". Novatec
#    Editora, 2008." run (same idea, spell-check run splitting).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "GUEDES, Gilleanes T. A. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "GUEDES, Gilleanes T. A. ",
    2) | Out-Null

$d.Content.Find.Execute(
    ". Novatec Editora, 2008.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Novatec Editora, 2008.",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) Merge split "Sistemas Conta Azul, Gran Money, EasyStore." run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Sistemas Conta Azul, Gran Money, EasyStore.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sistemas Conta Azul, Gran Money, EasyStore.",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) Add two new bullet paragraphs ("M0: mês atual." / "M-1: mês
#    anterior.") right after the "ES2: Disciplina de Engenharia de
#    Software 2." bullet and before the "Referências" heading.
# ------------------------------------------------------------------
$cnt = $d.Paragraphs.Count
$targetIdx = -1
for ($i = 1; $i -le $cnt; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*ES2*Disciplina de Engenharia de Software*") {
        $targetIdx = $i
    }
}

if ($targetIdx -gt 0) {
    $p = $d.Paragraphs($targetIdx)
    $p.Range.InsertParagraphAfter()

    $newPara1 = $d.Paragraphs($targetIdx + 1)
    $r1 = $d.Range($newPara1.Range.Start, $newPara1.Range.End - 1)
    $r1.InsertAfter("M0: mês atual.")

    $newPara1b = $d.Paragraphs($targetIdx + 1)
    $newPara1b.Range.InsertParagraphAfter()

    $newPara2 = $d.Paragraphs($targetIdx + 2)
    $r2 = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)
    $r2.InsertAfter("M-1: mês anterior.")
}
